$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4 (ALC)
$ws.Range("H4").Value = 4200.25
$ws.Range("I4").Value = 3600.3333
$ws.Range("K4").Value = 3600.3333
$ws.Range("M4").Value = -3486.3333

# Row 5 (ALC)
$ws.Range("H5").Value = 89.333336
$ws.Range("I5").Value = 89.333336
$ws.Range("K5").Value = 89.333336
$ws.Range("M5").Value = 25.666664

# Row 15 (ALC)
$ws.Range("H15").Value = 1772.4043
$ws.Range("I15").Value = 1772.4043
$ws.Range("K15").Value = 5317.2129
$ws.Range("M15").Value = -5148.2129

# Row 17 (ALC)
$ws.Range("H17").Value = 4994.1665
$ws.Range("J17").Value = 4994.1665
$ws.Range("L17").Value = 14982.4995
$ws.Range("N17").Value = -15318.4995

# Row 39 (ALC)
$ws.Range("H39").Value = 497.7143
$ws.Range("I39").Value = 534.8461
$ws.Range("J39").Value = 15
$ws.Range("K39").Value = 1604.5383
$ws.Range("L39").Value = 45
$ws.Range("M39").Value = -1308.5383
$ws.Range("N39").Value = -637

# Row 112 (ALC)
$ws.Range("H112").Value = 1016.2857
$ws.Range("I112").Value = 1019
$ws.Range("J112").Value = 1015.8333
$ws.Range("K112").Value = 3057
$ws.Range("L112").Value = 3047.4999
$ws.Range("N112").Value = -5263.4999
$ws.Range("M112").Value = -1949

# Row 138 (ALC)
$ws.Range("H138").Value = 3429.8765
$ws.Range("I138").Value = 3813.1428
$ws.Range("J138").Value = 3349.791
$ws.Range("K138").Value = 11439.4284
$ws.Range("L138").Value = 10049.373
$ws.Range("M138").Value = -6299.428400000001
$ws.Range("N138").Value = -20329.373

# Row 141 (ALC)
$ws.Range("H141").Value = 3392.9092
$ws.Range("I141").Value = 3057.2
$ws.Range("K141").Value = 9171.599999999999
$ws.Range("M141").Value = -3991.599999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1504.2916
$ws.Range("I2").Value = 727.6667
$ws.Range("K2").Value = 727.6667
$ws.Range("M2").Value = -614.6667

# Row 32 (ARM)
$ws.Range("H32").Value = 4380.147
$ws.Range("I32").Value = 3903.9375
$ws.Range("K32").Value = 3903.9375
$ws.Range("M32").Value = -3616.9375

# Row 45 (ARM)
$ws.Range("H45").Value = 1872.4546
$ws.Range("J45").Value = 1999.5
$ws.Range("L45").Value = 1999.5
$ws.Range("N45").Value = -2753.5

# Row 74 (ARM)
$ws.Range("H74").Value = 2878.7856
$ws.Range("I74").Value = 2726.4546
$ws.Range("J74").Value = 3437.3333
$ws.Range("K74").Value = 2726.4546
$ws.Range("L74").Value = 3437.3333
$ws.Range("M74").Value = -1852.4546
$ws.Range("N74").Value = -5185.3333

# Row 77 (ARM)
$ws.Range("H77").Value = 2878.7856
$ws.Range("I77").Value = 2726.4546
$ws.Range("J77").Value = 3437.3333
$ws.Range("K77").Value = 13632.273
$ws.Range("L77").Value = 17186.6665
$ws.Range("M77").Value = -9264.273000000001
$ws.Range("N77").Value = -25922.6665

# Row 88 (ARM)
$ws.Range("H88").Value = 1234.75
$ws.Range("I88").Value = 703.8
$ws.Range("K88").Value = 703.8
$ws.Range("M88").Value = -297.8

# Row 91 (ARM)
$ws.Range("H91").Value = 1234.75
$ws.Range("I91").Value = 703.8
$ws.Range("K91").Value = 703.8
$ws.Range("M91").Value = 700.2

# Row 116 (ARM)
$ws.Range("H116").Value = 1504.2916
$ws.Range("I116").Value = 727.6667
$ws.Range("K116").Value = 727.6667
$ws.Range("M116").Value = 1566.3333

# Row 132 (ARM)
$ws.Range("H132").Value = 1978.0834
$ws.Range("I132").Value = 1846.0869
$ws.Range("K132").Value = 5538.2607
$ws.Range("M132").Value = -3008.2607

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1504.2916
$ws.Range("I3").Value = 727.6667
$ws.Range("K3").Value = 727.6667
$ws.Range("M3").Value = -613.6667

# Row 94 (BSM)
$ws.Range("H94").Value = 600
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Row 107 (BSM)
$ws.Range("H107").Value = 3299.2856
$ws.Range("I107").Value = 3031.6667
$ws.Range("J107").Value = 4905
$ws.Range("K107").Value = 3031.6667
$ws.Range("L107").Value = 4905
$ws.Range("M107").Value = -1111.6667
$ws.Range("N107").Value = -8745

# Row 134 (BSM)
$ws.Range("H134").Value = 1641.3334
$ws.Range("I134").Value = 1749
$ws.Range("K134").Value = 5247
$ws.Range("M134").Value = -2712

$ws = $wb.Worksheets.Item("CRP")
# Row 2 (CRP)
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# Row 22 (CRP)
$ws.Range("H22").Value = 42653.168
$ws.Range("I22").Value = 185
$ws.Range("J22").Value = 63887.25
$ws.Range("K22").Value = 185
$ws.Range("L22").Value = 63887.25
$ws.Range("M22").Value = 165
$ws.Range("N22").Value = -64587.25

# Row 29 (CRP)
$ws.Range("H29").Value = 1519.6
$ws.Range("J29").Value = 1649.75
$ws.Range("L29").Value = 1649.75
$ws.Range("N29").Value = -2235.75

# Row 31 (CRP)
$ws.Range("H31").Value = 3478.5625
$ws.Range("I31").Value = 3216.7144
$ws.Range("K31").Value = 3216.7144
$ws.Range("M31").Value = -2921.7144

# Row 34 (CRP)
$ws.Range("H34").Value = 3478.5625
$ws.Range("I34").Value = 3216.7144
$ws.Range("K34").Value = 3216.7144
$ws.Range("M34").Value = -3014.7144

# Row 62 (CRP)
$ws.Range("H62").Value = 2399.5
$ws.Range("I62").Value = 2399.5
$ws.Range("K62").Value = 2399.5
$ws.Range("M62").Value = -1775.5

# Row 65 (CRP)
$ws.Range("H65").Value = 2399.5
$ws.Range("I65").Value = 2399.5
$ws.Range("K65").Value = 11997.5
$ws.Range("M65").Value = -8877.5

# Row 68 (CRP)
$ws.Range("H68").Value = 47500
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71 (CRP)
$ws.Range("H71").Value = 47500
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 134 (CRP)
$ws.Range("H134").Value = 2056.138
$ws.Range("I134").Value = 2003.28
$ws.Range("K134").Value = 6009.84
$ws.Range("M134").Value = -3474.84

$ws = $wb.Worksheets.Item("CUL")
# Row 13 (CUL)
$ws.Range("H13").Value = 3000
$ws.Range("I13").Value = 3000
$ws.Range("K13").Value = 9000
$ws.Range("M13").Value = -8832

# Row 97 (CUL)
$ws.Range("H97").Value = 6945754.5
$ws.Range("J97").Value = 8930092
$ws.Range("L97").Value = 26790276
$ws.Range("N97").Value = -26791268

# Row 98 (CUL)
$ws.Range("H98").Value = 845.3333
$ws.Range("I98").Value = 250
$ws.Range("J98").Value = 1143
$ws.Range("K98").Value = 750
$ws.Range("L98").Value = 3429
$ws.Range("M98").Value = 748
$ws.Range("N98").Value = -6425

# Row 131 (CUL)
$ws.Range("H131").Value = 2639
$ws.Range("I131").Value = 1873.75
$ws.Range("J131").Value = 2917.2727
$ws.Range("K131").Value = 5621.25
$ws.Range("L131").Value = 8751.8181
$ws.Range("M131").Value = -581.25
$ws.Range("N131").Value = -18831.8181

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (GSM)
$ws.Range("H113").Value = 1491.5
$ws.Range("J113").Value = 999
$ws.Range("L113").Value = 999
$ws.Range("N113").Value = -5339

# Row 132 (GSM)
$ws.Range("H132").Value = 1159.5834
$ws.Range("I132").Value = 1159.5834
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3478.7502
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -948.7501999999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 23 (LTW)
$ws.Range("H23").Value = 93333
$ws.Range("I23").Value = 93333
$ws.Range("K23").Value = 93333
$ws.Range("M23").Value = -93103

# Row 29 (LTW)
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 100 (LTW)
$ws.Range("H100").Value = 4374.7144
$ws.Range("I100").Value = 2874.3333
$ws.Range("K100").Value = 2874.3333
$ws.Range("M100").Value = -2333.3333

# Row 132 (LTW)
$ws.Range("H132").Value = 3287.9167
$ws.Range("I132").Value = 2295.1
$ws.Range("J132").Value = 8252
$ws.Range("K132").Value = 6885.299999999999
$ws.Range("L132").Value = 24756
$ws.Range("M132").Value = -4355.299999999999
$ws.Range("N132").Value = -29816

$ws = $wb.Worksheets.Item("WVR")
# Row 34 (WVR)
$ws.Range("H34").Value = 26666
$ws.Range("I34").Value = 26666
$ws.Range("K34").Value = 26666
$ws.Range("M34").Value = -26463

# Row 132 (WVR)
$ws.Range("H132").Value = 2090.9443
$ws.Range("I132").Value = 2265.4443
$ws.Range("J132").Value = 1916.4445
$ws.Range("K132").Value = 6796.3329
$ws.Range("L132").Value = 5749.333500000001
$ws.Range("M132").Value = -4266.3329
$ws.Range("N132").Value = -10809.3335

# Row 136 (WVR)
$ws.Range("H136").Value = 4373.722
$ws.Range("I136").Value = 4395.5
$ws.Range("J136").Value = 4199.5
$ws.Range("K136").Value = 13186.5
$ws.Range("L136").Value = 12598.5
$ws.Range("M136").Value = -10636.5
$ws.Range("N136").Value = -17698.5
